# Daily attendance processing
# For every row in the "Recorded By" column (G), reorder the comma-separated
# list of recorder names so that any "System"/"system" entries move to the
# front of the list, while keeping the relative order of all entries
# (System entries first, then the remaining entries in their original order).
# Rows whose "Recorded By" value has no "System" entry (single value, or a
# list without "System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val.Split(",")
    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemParts += $trimmed
        } else {
            $otherParts += $trimmed
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $newParts = $systemParts + $otherParts
    $newVal = $newParts -join ", "
    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
